$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    at the top of the data (row 2), shifting everything else down,
#    and bump the running index in column A for the shifted rows.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The row-insert can bleed formatting in from neighbouring rows; reset
# the plain data cells (B2:D2) back to the default/no style first.
$summary.Range("B2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.46

# Give A2 the same look (bold + border) as the rest of column A by
# copying the format from A3 (which already carries it after the shift).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$summary.Range("A2").Value = 0

# Bump the running index (column A) of all the rows that shifted down.
for ($r = 3; $r -le 8; $r++) {
    $cell = $summary.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# ------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q3" right before "2022-Q2"
#    and populate it with the quarterly fund-holder table.
# ------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

# Re-fetch a fresh handle by name (worksheet references can go stale
# once the sheet collection is mutated).
$q3 = $wb.Worksheets.Item("2022-Q3")

# Copy the bold/bordered header style used throughout the workbook
# (sourced from the "总计" sheet's header row / index column) onto the
# header row and the index column (A) of the new sheet.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows. Columns D, E, F, G (and the fund code in B) are stored as
# plain text in the source workbook even though they look numeric, so
# they are written with a leading apostrophe to force text and then
# the style is reset to "Normal" so no stray number-format sticks.
$codes  = @("160322","012884","002332","007291","006537","003993","002333","004321","006477")
$names  = @("华夏港股通精选股票（LOF）","华夏港股通精选股票(LOF) C","汇丰晋信沪港深股票A","汇丰晋信港股通双核策略混合","恒生前海港股通精选混合","前海开源沪港深核心驱动灵活配置混合","汇丰晋信沪港深股票C","前海开源沪港深强国产业灵活配置混合","中邮沪港深精选混合")
$scale  = @("13.65","13.65","5.45","6.17","0.95","0.53","0.48","0.11","0.06")
$pos    = @("84.60","84.60","94.22","94.56","90.50","82.41","94.22","78.52","90.21")
$ratio  = @("2.97","2.97","5.09","4.41","3.64","5.86","5.09","5.54","5.62")
$value  = @("0.4054","0.4054","0.2774","0.2721","0.0346","0.0311","0.0244","0.0061","0.0034")
$rank   = @(7,7,6,9,9,10,6,4,7)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2

    $q3.Cells.Item($r, 1).Value = $i

    $q3.Cells.Item($r, 2).Value = "'" + $codes[$i]
    $q3.Cells.Item($r, 2).Style = "Normal"

    $q3.Cells.Item($r, 3).Value = $names[$i]

    $q3.Cells.Item($r, 4).Value = "'" + $scale[$i]
    $q3.Cells.Item($r, 4).Style = "Normal"

    $q3.Cells.Item($r, 5).Value = "'" + $pos[$i]
    $q3.Cells.Item($r, 5).Style = "Normal"

    $q3.Cells.Item($r, 6).Value = "'" + $ratio[$i]
    $q3.Cells.Item($r, 6).Style = "Normal"

    $q3.Cells.Item($r, 7).Value = "'" + $value[$i]
    $q3.Cells.Item($r, 7).Style = "Normal"

    $q3.Cells.Item($r, 8).Value = $rank[$i]
}

Write-Output "edit complete"
